# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.553.33"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.504.98"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'575.19"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'166.72"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "2.505.28"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'0.360"
$ws.Range("E12").Value = "  +5.51%  "
$ws.Range("D13").Value = "'4.92"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "2.965.55"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "69.421.73"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "'24.82"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "2.506.70"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'11.26"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'7.57"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "'350.23"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'70.44"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "'8.83"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").Value = "2.638.62"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.0₃0890"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'460.39"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'159.72"
$ws.Range("E36").Value = "  +4.19%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'18.44"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'4.68"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "'38.22"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("D46").Value = "'1.10"
$ws.Range("E46").Value = "  -7.10%  "
$ws.Range("D47").Value = "'142.18"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "'0.520"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("D50").Value = "'0.0734"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +2.86%  "
